$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet updates ---

# Row 2
$wsForecast.Range("L2").Value = 0.95

# Row 3
$wsForecast.Range("D3").Value = 54
$wsForecast.Range("H3").Value = 9.93
$wsForecast.Range("L3").Value = 1.15

# Row 4
$wsForecast.Range("H4").Value = 10.31
$wsForecast.Range("L4").Value = 0.87

# Row 5
$wsForecast.Range("H5").Value = 8.369999999999999
$wsForecast.Range("L5").Value = 0.87

# Row 6
$wsForecast.Range("H6").Value = 7.56
$wsForecast.Range("L6").Value = 1.15

# Row 7
$wsForecast.Range("H7").Value = 6.74
$wsForecast.Range("L7").Value = 0.85

# Row 8
$wsForecast.Range("H8").Value = 5.74
$wsForecast.Range("L8").Value = 0.83

# Row 9
$wsForecast.Range("D9").Value = 45
$wsForecast.Range("H9").Value = 5.2
$wsForecast.Range("L9").Value = 0.85

# Row 10
$wsForecast.Range("H10").Value = 4.04
$wsForecast.Range("L10").Value = 1.09

# Row 11
$wsForecast.Range("H11").Value = 3.04
$wsForecast.Range("L11").Value = 0.97

# Row 12
$wsForecast.Range("H12").Value = 2.04
$wsForecast.Range("L12").Value = 0.88

# Row 13
$wsForecast.Range("H13").Value = 1.01
$wsForecast.Range("J13").Value = "Normal"
$wsForecast.Range("L13").Value = 0.86

# Row 14
$wsForecast.Range("H14").Value = 0.01
$wsForecast.Range("L14").Value = 0.96

# Row 15
$wsForecast.Range("D15").Value = 31
$wsForecast.Range("L15").Value = 1.01

# Row 16
$wsForecast.Range("D16").Value = 35
$wsForecast.Range("L16").Value = 1.12

# Row 17
$wsForecast.Range("L17").Value = 0.8100000000000001

# --- Summary sheet updates ---
# These cells store numeric-looking values as TEXT (matching the source
# workbook's inline-string cells), so force text formatting before
# assigning, otherwise Excel auto-coerces the numeric-looking string to a
# real number.
$summaryTextCells = @{
    "B9"  = "738"
    "B10" = "398"
    "B11" = "204"
    "B12" = "54"
    "B14" = "31"
}
foreach ($addr in $summaryTextCells.Keys) {
    $cell = $wsSummary.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryTextCells[$addr]
}
